$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# --- Experience section: rewrite the single summary sentence, mention the company ---
$oldExperience = "Developed an AI web-app for resume enhancement using LLM models, training them to identify required resume features and fine-tuning with targeted prompts for content generation, resulting in improved content creation efficiency."
$newExperience = "Developed An AI Web-App For Resume Enhancement Using LLM Models At Blue Silicon. Trained Models To Identify Required Resume Features And Fine-Tuned With Targeted Prompts For Content Generation, Resulting In Improved Content Creation Efficiency."
Replace-Text $oldExperience $newExperience

# --- Education section: combine the two line-broken entries into one sentence ---
$oldEducation = "Bachelor of Engineering in Computer Science, AVIT, May 2026, GPA 7.1;^lIntermediate, Narayana Jr College, June 2022, GPA 8.9"
$newEducation = "Bachelor Of Engineering In Computer Science Was Earned From Avit, Graduating May 2026 With A Gpa Of 7.1. Intermediate Education Was Completed At Narayana Jr College, Graduating June 2022 With A Gpa Of 8.9."
Replace-Text $oldEducation $newEducation

# --- Skills section: replace the multi-line bulleted list (line-break overlap bug) with one line ---
$oldSkills = "Good with people, can talk and listen well but sometimes nervous,^lKnow some billing and petty cash handling but not expert,^lExperience running group sessions for mental health and social skills,^lCan plan and do social activities but not very creative,^l^lAble to multitask but sometimes get overwhelmed,^l^lBasic computer skills like Microsoft Word and Excel, "
$newSkills = "Python 100%, Git (kinda), Linux? sort of, Web scraping I guess, Can use Stack Overflow good, Copy paste code really fast, Debugging (sometimes works), Not great at Java but maybe later"
Replace-Text $oldSkills $newSkills

# --- Projects section: combine all projects into one paragraph with double line breaks between them ---
$oldProjects = "Developed and implemented a QR scanner and generator, leveraging TypeScript and Node.js for the backend. Spearheaded backend development, ensuring seamless functionality. Additionally, contributed to Prediction Pro, a full-stack application built with TypeScript, React, and PostgreSQL, demonstrating expertise in modern technologies and collaborative skills through frontend and database integration. These projects showcased versatility in tech stacks, full-stack development capabilities, and effective collaboration."
$newProjects = "Developed A Snake Game Using Basic Programming Concepts.^l^lCreated A Weather App That Retrieves Data From An API; However, Its Functionality Was Limited When The API Stopped Working.^l^lDesigned A Discord Bot That Automates Responses And Interactions.^l^lWrote An Automation Script That Efficiently Renames Files Using Scripting Techniques.^l^lAttempted To Develop A Web Scraper, But Encountered Blocking Issues.^l^lBuilt A Calculator With Fundamental Arithmetic Operations."
Replace-Text $oldProjects $newProjects
